$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 20 new rows right after the header row, pushing existing data down
$ws.Rows("2:21").Insert(-4121)

# Force Text format on the new rows so phone numbers / DDDs / dates are
# stored as plain text (matching the rest of the sheet) instead of being
# auto-converted to numbers / dates by Excel
$newRowsRange = $ws.Range("A2:C21")
$newRowsRange.NumberFormat = "@"

# Populate the new rows with the newly added signups
$ws.Range("A2").Value = "+5513997050892"
$ws.Range("B2").Value = "13"
$ws.Range("C2").Value = "2024-10-18"
$ws.Range("A3").Value = "+5511977538423"
$ws.Range("B3").Value = "11"
$ws.Range("C3").Value = "2024-10-18"
$ws.Range("A4").Value = "+5511920143449"
$ws.Range("B4").Value = "11"
$ws.Range("C4").Value = "2024-10-18"
$ws.Range("A5").Value = "+556194520205"
$ws.Range("B5").Value = "61"
$ws.Range("C5").Value = "2024-10-18"
$ws.Range("A6").Value = "+5511947132973"
$ws.Range("B6").Value = "11"
$ws.Range("C6").Value = "2024-10-18"
$ws.Range("A7").Value = "+5521974855968"
$ws.Range("B7").Value = "21"
$ws.Range("C7").Value = "2024-10-18"
$ws.Range("A8").Value = "+556299060586"
$ws.Range("B8").Value = "62"
$ws.Range("C8").Value = "2024-10-18"
$ws.Range("A9").Value = "+555193669200"
$ws.Range("B9").Value = "51"
$ws.Range("C9").Value = "2024-10-18"
$ws.Range("A10").Value = "+5511948066627"
$ws.Range("B10").Value = "11"
$ws.Range("C10").Value = "2024-10-18"
$ws.Range("A11").Value = "+5517981320321"
$ws.Range("B11").Value = "17"
$ws.Range("C11").Value = "2024-10-18"
$ws.Range("A12").Value = "+5521994172524"
$ws.Range("B12").Value = "21"
$ws.Range("C12").Value = "2024-10-18"
$ws.Range("A13").Value = "+5519994757305"
$ws.Range("B13").Value = "19"
$ws.Range("C13").Value = "2024-10-18"
$ws.Range("A14").Value = "+5514996136195"
$ws.Range("B14").Value = "14"
$ws.Range("C14").Value = "2024-10-18"
$ws.Range("A15").Value = "+555199322869"
$ws.Range("B15").Value = "51"
$ws.Range("C15").Value = "2024-10-18"
$ws.Range("A16").Value = "+5519989071084"
$ws.Range("B16").Value = "19"
$ws.Range("C16").Value = "2024-10-18"
$ws.Range("A17").Value = "+558491181674"
$ws.Range("B17").Value = "84"
$ws.Range("C17").Value = "2024-10-18"
$ws.Range("A18").Value = "+5511994736173"
$ws.Range("B18").Value = "11"
$ws.Range("C18").Value = "2024-10-18"
$ws.Range("A19").Value = "+5513991481556"
$ws.Range("B19").Value = "13"
$ws.Range("C19").Value = "2024-10-18"
$ws.Range("A20").Value = "+5521983285325"
$ws.Range("B20").Value = "21"
$ws.Range("C20").Value = "2024-10-17"
$ws.Range("A21").Value = "+556195021826"
$ws.Range("B21").Value = "61"
$ws.Range("C21").Value = "2024-10-17"

# Re-apply the standard data-row style (as used by row 22, the first
# pre-existing data row) to the newly inserted rows so formatting matches
$templateRow = $ws.Range("A22:C22")
$templateRow.Copy()
$newRowsRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false
